$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.289.91"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.285.00"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.32"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.91"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.33"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.31"
$ws.Range("E12").Value = "  +7.59%  "
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.50"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.627.05"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.867"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.277.23"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.299.94"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000109"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.81"
$ws.Range("E20").Value = "  +4.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.09"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.46"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.04"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.59"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.86"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.44"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.10"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.36"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.83"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.62"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0904"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  +3.94%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.128"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.64"
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0369"
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.63"
$ws.Range("E40").Value = "  +8.97%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.12"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.16"
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.240"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.16"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.38"
$ws.Range("E46").Value = "  -4.66%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.70"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.60"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.601"
$ws.Range("E51").Value = "  +9.85%  "
